# Turn the old "Heading1 + bold byline" header into a pandoc-style title
# block: a Title-styled heading paragraph and an Authors-styled byline
# paragraph, each split word-by-word into separate runs (as the target
# OOXML does). The two legacy <w:bookmarkStart/End> markers that used to
# wrap the heading paragraph are also removed.

$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function New-WordRunsXml($words) {
    $xml = ""
    foreach ($word in $words) {
        $escaped = $word -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        $xml += "<w:r><w:t xml:space='preserve'>$escaped</w:t></w:r>"
    }
    return $xml
}

# --- Paragraph 1: "On Pilgrimage - October 1963" (Heading1, bookmarked) ---
# Delete the paragraph's whole text+mark first, which leaves the two
# bookmark markers that used to wrap it collapsed at document position 0
# with no paragraph around them any more.
$d.Paragraphs(1).Range.Delete()

# A zero-length delete sitting exactly at position 0 strips whichever
# bookmark marker is anchored there without touching any text; do it
# twice to clear both the start and end markers.
$d.Range(0, 0).Delete(1, 1)
$d.Range(0, 0).Delete(1, 1)

# Insert the new Title paragraph, word-by-word as separate runs, at the
# (now bookmark-free) start of the document.
$titleWords = @("On", " ", "Pilgrimage", " ", "-", " ", "October", " ", "1963")
$titleXml = "<w:p $wns><w:pPr><w:pStyle w:val='Title'/></w:pPr>" + (New-WordRunsXml $titleWords) + "</w:p>"
$d.Range(0, 0).InsertXML($titleXml)

# --- Paragraph 2: "By Dorothy Day" (bold) -> Authors paragraph "Dorothy Day" ---
$authorWords = @("Dorothy", " ", "Day")
$authorXml = "<w:p $wns><w:pPr><w:pStyle w:val='Authors'/></w:pPr>" + (New-WordRunsXml $authorWords) + "</w:p>"
$d.Paragraphs(2).Range.InsertXML($authorXml)
